$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data for user view review row
$ws.Range("E8").Value = 43013
$ws.Range("E9").Value = 43013

$ws.Range("C10").Value = "user view review"
$ws.Range("D10").Value = "kleine Nachbesserungen"
$ws.Range("E10").Value = 43026

$ws.Range("E7").Copy()
$ws.Range("E8:E10").PasteSpecial(-4122)

$ws.Range("C10:D10").Style = $ws.Range("C9:D9").Style

$ws.Range("D16").Select()
